# Applies updated "My 4 Weeks Forecast" (column C) values per the
# "added 4wk low sales check" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Value = 3 },
    @{ Row = 3; Value = 212.6 },
    @{ Row = 4; Value = 4 },
    @{ Row = 5; Value = 2.6 },
    @{ Row = 6; Value = 2 },
    @{ Row = 8; Value = 516 },
    @{ Row = 9; Value = 569 },
    @{ Row = 10; Value = 36 },
    @{ Row = 11; Value = 527 },
    @{ Row = 12; Value = 3.800000000000001 },
    @{ Row = 13; Value = 6.4 },
    @{ Row = 14; Value = 2542 },
    @{ Row = 15; Value = 140 },
    @{ Row = 16; Value = 1.6 },
    @{ Row = 18; Value = 9.6 },
    @{ Row = 19; Value = 208 },
    @{ Row = 20; Value = 2 },
    @{ Row = 21; Value = 498 },
    @{ Row = 22; Value = 2 },
    @{ Row = 23; Value = 250.6 },
    @{ Row = 25; Value = 0 },
    @{ Row = 26; Value = 800 },
    @{ Row = 27; Value = 0 },
    @{ Row = 28; Value = 362 },
    @{ Row = 29; Value = 93 },
    @{ Row = 30; Value = 34 },
    @{ Row = 31; Value = 4 },
    @{ Row = 32; Value = 662 },
    @{ Row = 33; Value = 545 },
    @{ Row = 34; Value = 1.2 },
    @{ Row = 38; Value = 53 },
    @{ Row = 39; Value = 1227 },
    @{ Row = 40; Value = 0 },
    @{ Row = 41; Value = 3.4 },
    @{ Row = 43; Value = 66 },
    @{ Row = 44; Value = 22.8 },
    @{ Row = 45; Value = 12.6 },
    @{ Row = 46; Value = 359 },
    @{ Row = 47; Value = 0.8 },
    @{ Row = 48; Value = 3.6 },
    @{ Row = 49; Value = 4 },
    @{ Row = 50; Value = 832 },
    @{ Row = 51; Value = 182 },
    @{ Row = 52; Value = 2 },
    @{ Row = 53; Value = 153 },
    @{ Row = 54; Value = 115 },
    @{ Row = 55; Value = 0.8 },
    @{ Row = 56; Value = 17.2 },
    @{ Row = 57; Value = 0 },
    @{ Row = 58; Value = 237.4 },
    @{ Row = 59; Value = 76 },
    @{ Row = 60; Value = 57 },
    @{ Row = 61; Value = 0 },
    @{ Row = 62; Value = 388 },
    @{ Row = 63; Value = 6 },
    @{ Row = 64; Value = 57 },
    @{ Row = 65; Value = 14.8 },
    @{ Row = 66; Value = 96 },
    @{ Row = 67; Value = 128 },
    @{ Row = 69; Value = 218.8 },
    @{ Row = 70; Value = 110 },
    @{ Row = 71; Value = 167 },
    @{ Row = 72; Value = 353 },
    @{ Row = 73; Value = 136 },
    @{ Row = 74; Value = 5.200000000000001 },
    @{ Row = 75; Value = 53.2 },
    @{ Row = 76; Value = 92 },
    @{ Row = 77; Value = 622.5999999999999 },
    @{ Row = 79; Value = 363 },
    @{ Row = 80; Value = 55.40000000000001 },
    @{ Row = 81; Value = 30.8 },
    @{ Row = 82; Value = 340 },
    @{ Row = 83; Value = 1627 },
    @{ Row = 84; Value = 12.8 },
    @{ Row = 85; Value = 9.199999999999999 },
    @{ Row = 86; Value = 7.2 },
    @{ Row = 87; Value = 181.6 },
    @{ Row = 88; Value = 1089.8 },
    @{ Row = 89; Value = 4 },
    @{ Row = 90; Value = 28 },
    @{ Row = 91; Value = 448 },
    @{ Row = 92; Value = 8.4 },
    @{ Row = 93; Value = 6 },
    @{ Row = 94; Value = 675.4000000000001 },
    @{ Row = 95; Value = 186 },
    @{ Row = 96; Value = 51 },
    @{ Row = 97; Value = 1605.4 },
    @{ Row = 98; Value = 1199.2 },
    @{ Row = 99; Value = 4512.200000000001 },
    @{ Row = 100; Value = 342 },
    @{ Row = 101; Value = 413.8 },
    @{ Row = 102; Value = 799.8 },
    @{ Row = 103; Value = 313 },
    @{ Row = 104; Value = 146.4 },
    @{ Row = 105; Value = 1863.2 },
    @{ Row = 106; Value = 103 },
    @{ Row = 107; Value = 4 },
    @{ Row = 108; Value = 171.2 },
    @{ Row = 109; Value = 7 },
    @{ Row = 110; Value = 251 },
    @{ Row = 111; Value = 245 },
    @{ Row = 112; Value = 100.2 },
    @{ Row = 113; Value = 39.2 },
    @{ Row = 114; Value = 21.6 },
    @{ Row = 115; Value = 24 },
    @{ Row = 116; Value = 24 },
    @{ Row = 117; Value = 52 },
    @{ Row = 118; Value = 10.4 },
    @{ Row = 119; Value = 48 },
    @{ Row = 120; Value = 17.6 },
    @{ Row = 121; Value = 38.40000000000001 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Value
}
